# LSh enabled on prod
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Turn off concurrent calculation (workbook calcPr concurrentCalc="0")
$excel.Application.MultiThreadedCalculation.Enabled = $false

# Row 49: GX 460 Premium -> 2020 -> 2021, MSRP 80010 -> 90500
$ws.Range("C49").Value = 2021
$ws.Range("D49").Value = 90500

# Row 50: GX 460 Luxury -> 2020 -> 2021, MSRP 83230 -> 93750
$ws.Range("C50").Value = 2021
$ws.Range("D50").Value = 93750

# Rows 73-78: small MSRP bumps
$ws.Range("D73").Value = 33000
$ws.Range("D74").Value = 37700
$ws.Range("D75").Value = 35200
$ws.Range("D76").Value = 39900
$ws.Range("D77").Value = 35000
$ws.Range("D78").Value = 37200

# Row 80: MSRP bump
$ws.Range("D80").Value = 119900

# New row 100 (assign shared strings in the same order the original
# workbook introduces them: MODEL TRIM text first, then TRIM CODE)
$ws.Range("B100").Value = "LC 500 INSPIRATION SERIES"
$ws.Range("A100").Value = "9260LE"
$ws.Range("C100").Value = 2021
$ws.Range("D100").Value = 110420
$ws.Range("D100").NumberFormat = $ws.Range("D79").NumberFormat
$ws.Range("E100").Value = 1025
$ws.Range("E100").NumberFormat = $ws.Range("E99").NumberFormat

# Update sheet view to reflect the final selection/scroll position
$ws.Range("D73:D78").Select()
$excel.ActiveWindow.ScrollRow = 46
